# Update "想去人数" (F column) counts across the workbook sheets to reflect
# the latest scrape results output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 2651
    8  = 1225
    9  = 570
    10 = 309
    12 = 126
    14 = 5747
    15 = 1783
    16 = 4180
    20 = 4871
    21 = 6257
    25 = 3784
    26 = 499
    32 = 475
    33 = 562
    35 = 201
    36 = 1726
    39 = 1144
    43 = 3426
    48 = 18
    49 = 3891
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "本地生活" (local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 3941

# --- Sheet "全部类型" (all types, combined view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 3941
    3  = 2651
    12 = 1225
    13 = 570
    14 = 309
    15 = 126
    18 = 1783
    19 = 4180
    20 = 4871
    24 = 3784
    25 = 499
    31 = 475
    32 = 562
    35 = 201
    36 = 1726
    43 = 3426
    49 = 3891
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
